$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Jogo"
$headerVals = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15)
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headerVals[$i]
}

# Row 2
$row2Vals = @(2725,2,4,5,6,7,10,11,12,15,16,18,20,21,22,23)
for ($i = 0; $i -lt $row2Vals.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2Vals[$i]
}

# Row 3
$row3Vals = @(2726,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17)
for ($i = 0; $i -lt $row3Vals.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3Vals[$i]
}

# Clear old leftover cells beyond row 3 / old data (row4 baz/qux, and old D column data)
$ws.Range("A4:P4").ClearContents()

# The old C3 cell carried a date-number-format style; clear formatting so it
# reverts to the default style like the rest of the refreshed table.
$ws.Range("A1:P4").ClearFormats()

# Row heights as left by Excel after editing this table.
$ws.Range("A1:P4").RowHeight = 15.75

# Page margins matching the saved worksheet page setup.
$ws.PageSetup.LeftMargin = 0.511811024 * 72
$ws.PageSetup.RightMargin = 0.511811024 * 72
$ws.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$ws.PageSetup.FooterMargin = 0.31496062000000002 * 72

# Selection left on the sheet after the edits.
$ws.Range("Q3").Select()
